# Add two new worksheets "area_lores_basic" and "area_pop_sum_basic" at the
# end of the workbook, re-running the summary-stats calcs with updated
# (basic) geounits data, mirroring the structure/formatting of the existing
# "area_lores" and "area_pop_sum" sheets.

$wb = $excel.ActiveWorkbook

$wsAreaLores  = $wb.Worksheets.Item("area_lores")
$wsAreaPopSum = $wb.Worksheets.Item("area_pop_sum")

# --- New sheet 1: area_lores_basic (cloned layout of area_lores) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLoresBasic = $wb.Worksheets.Add($null, $lastSheet)
$wsLoresBasic.Name = "area_lores_basic"

$wsLoresBasic.Range("A1").Value = "index"
$wsLoresBasic.Range("B1").Value = "area"
$wsLoresBasic.Range("A2").Value = "count"
$wsLoresBasic.Range("B2").Value = 191
$wsLoresBasic.Range("A3").Value = "mean"
$wsLoresBasic.Range("B3").Value = 4.690043535632206
$wsLoresBasic.Range("A4").Value = "std"
$wsLoresBasic.Range("B4").Value = 5.404762470031221
$wsLoresBasic.Range("A5").Value = "min"
$wsLoresBasic.Range("B5").Value = 0.4491082895602994

# these labels look like percentages, so force them to stay as plain text
# and then strip the resulting number-format styling back off again
$wsLoresBasic.Range("A6:A8").NumberFormat = "@"
$wsLoresBasic.Range("A6").Value = "25%"
$wsLoresBasic.Range("B6").Value = 1.567172533707438
$wsLoresBasic.Range("A7").Value = "50%"
$wsLoresBasic.Range("B7").Value = 3.031139052727144
$wsLoresBasic.Range("A8").Value = "75%"
$wsLoresBasic.Range("B8").Value = 5.911424536988898
$wsLoresBasic.Range("A6:A8").ClearFormats()

$wsLoresBasic.Range("A9").Value = "max"
$wsLoresBasic.Range("B9").Value = 37.62246804805788

# copy header formatting (bold, bordered, centered) from area_lores!A1:B1
$wsAreaLores.Range("A1:B1").Copy()
$wsLoresBasic.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New sheet 2: area_pop_sum_basic (cloned layout of area_pop_sum) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPopSumBasic = $wb.Worksheets.Add($null, $lastSheet)
$wsPopSumBasic.Name = "area_pop_sum_basic"

$wsPopSumBasic.Range("A1").Value = "index"
$wsPopSumBasic.Range("B1").Value = 0
$wsPopSumBasic.Range("A2").Value = "area"
$wsPopSumBasic.Range("B2").Value = 895.7983153057514
$wsPopSumBasic.Range("A3").Value = "population"
$wsPopSumBasic.Range("B3").Value = 3303002
$wsPopSumBasic.Range("A4").Value = "density"
$wsPopSumBasic.Range("B4").Value = 3687.216132877665

# copy header formatting (bold, bordered, centered) from area_pop_sum!A1:B1
$wsAreaPopSum.Range("A1:B1").Copy()
$wsPopSumBasic.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
